$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns affected by the permutation: H, I, K, L, M, N, O, P, Q, R, S, T, U, V
$cols = @("H","I","K","L","M","N","O","P","Q","R","S","T","U","V")

# Snapshot the "before" values for each row/col into a hashtable keyed "row_col"
$snap = @{}
for ($r = 2; $r -le 41; $r++) {
    foreach ($c in $cols) {
        $snap["$r" + "_" + $c] = $ws.Range($c + $r).Value2
    }
}

# Destination row -> source row mapping (permutation), derived from the commit diff
$map = @{}
$map[2] = 15
$map[3] = 41
$map[4] = 40
$map[5] = 23
$map[6] = 21
$map[7] = 2
$map[8] = 11
$map[9] = 8
$map[10] = 31
$map[11] = 38
$map[12] = 4
$map[13] = 14
$map[14] = 29
$map[15] = 32
$map[16] = 20
$map[17] = 3
$map[18] = 35
$map[19] = 39
$map[20] = 26
$map[21] = 36
$map[22] = 22
$map[23] = 16
$map[24] = 27
$map[25] = 25
$map[26] = 19
$map[27] = 37
$map[28] = 34
$map[29] = 6
$map[30] = 13
$map[31] = 10
$map[32] = 12
$map[33] = 30
$map[34] = 5
$map[35] = 28
$map[36] = 7
$map[37] = 33
$map[38] = 24
$map[39] = 9
$map[40] = 17
$map[41] = 18

# Apply: each destination row gets the snapshot values of its mapped source row
for ($r = 2; $r -le 41; $r++) {
    $src = $map[$r]
    foreach ($c in $cols) {
        $ws.Range($c + $r).Value2 = $snap["$src" + "_" + $c]
    }
}

